$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Row 18 (Min): B column now uses MIN instead of MAX ---
$ws.Range("B18").Formula = "=MIN(timeBuildBatch!`$H:`$H)"

# --- Row 19 (Average): B column now uses AVERAGE instead of MAX ---
$ws.Range("B19").Formula = "=AVERAGE(timeBuildBatch!`$H:`$H)"

# --- Row 20 (StdDev(p)): B column now uses STDEV.P instead of MAX ---
$ws.Range("B20").Formula = "=STDEV.P(timeBuildBatch!`$H:`$H)"

# --- New rows 23/24: total sec / total days computed from the grand total and average rate ---
$ws.Range("A23").Value = "total sec"
$ws.Range("B23").Formula = "=`$C`$14/B19"
$ws.Range("B23").Style = "Normal"
$ws.Range("D23").Value = "sec/day"
$ws.Range("E23").Value = 86400

$ws.Range("A24").Value = "total days"
$ws.Range("B24").Formula = "=B23/E23"
$ws.Range("B24").Style = "Normal"

# --- G14: add the "dubious" note next to the File rates (files / sec) total ---
$ws.Range("G14").Value = "?? Dubious. Something's off"

# --- Column widths: widen C (now holds longer labels) and give E a width too ---
$ws.Columns("C").ColumnWidth = 20.8
$ws.Columns("E").ColumnWidth = 13.3

# --- Leave the selection on the last cell touched, like the author did ---
$ws.Range("C24").Select()
